$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cost Center")

# Replace the Cost Center codes (CC001..CC006) with the new Automated Test codes
$ws.Range("C2").Value = "AutomatedTest001"
$ws.Range("C3").Value = "AutomatedTest002"
$ws.Range("C4").Value = "AutomatedTest003"
$ws.Range("C5").Value = "AutomatedTest004"
$ws.Range("C6").Value = "AutomatedTest005"
$ws.Range("C7").Value = "AutomatedTest006"

# Update the entity codes for rows 6 and 7
$ws.Range("B6").Value = "INAF"
$ws.Range("B7").Value = "KAEF"

# Move the active selection to D11
$ws.Range("D11").Select() | Out-Null
